$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2: shorten "La Mojana-Camellones" -> "Camellones"
$ws.Range("D2").Value = "Camellones"

# Fix typo Cordoba -> Córdoba in the description cell (H2)
$ws.Range("H2").Value = 'Shapefile, geometria tipo línea con información de la localización de los "camellones", una técnica agrícola ancestral que favorece la productividad en entornos de manejo de agua, así como los predios categorizados por su "área de influencia media" y "directa"  de los departamentos de Córdoba (municipio de Ayapel)  y Sucre (Majagual, San Benito Abad, San Marcos y Sucre)'

# Enable wrap text on the description cell (H2) to match new style
$ws.Range("H2").WrapText = $true

# Move the active selection from P2 to I2
$ws.Range("I2").Select()
